$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 782-783, shifting existing rows 782-860 down to 784-862
$ws.Rows("782:783").Insert()

# Row 782 - new data
$ws.Range("A782").Value = 4
$ws.Range("B782").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C782").Value = 'Los Lagos'
$ws.Range("D782").Value = 45166
$ws.Range("E782").Value = 10
$ws.Range("F782").Value = 'Fruta'
$ws.Range("G782").Value = 100102
$ws.Range("H782").Value = 'Cítricos'
$ws.Range("I782").Value = 100102005
$ws.Range("J782").Value = 'Naranja'
$ws.Range("K782").Value = 'Navel Late'
$ws.Range("L782").Value = 'Primera'
$ws.Range("M782").Value = 200
$ws.Range("N782").Value = 18000
$ws.Range("O782").Value = 18000
$ws.Range("P782").Value = 18000
$ws.Range("Q782").Value = '$/caja 15 kilos empedrada'
$ws.Range("R782").Value = 'Región de O''Higgins'
$ws.Range("S782").Value = 1200
$ws.Range("T782").Value = 15

# Row 783 - new data
$ws.Range("A783").Value = 4
$ws.Range("B783").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C783").Value = 'Los Lagos'
$ws.Range("D783").Value = 45166
$ws.Range("E783").Value = 10
$ws.Range("F783").Value = 'Fruta'
$ws.Range("G783").Value = 100102
$ws.Range("H783").Value = 'Cítricos'
$ws.Range("I783").Value = 100102005
$ws.Range("J783").Value = 'Naranja'
$ws.Range("K783").Value = 'Navel Late'
$ws.Range("L783").Value = 'Segunda'
$ws.Range("M783").Value = 200
$ws.Range("N783").Value = 15000
$ws.Range("O783").Value = 15000
$ws.Range("P783").Value = 15000
$ws.Range("Q783").Value = '$/caja 15 kilos empedrada'
$ws.Range("R783").Value = 'Región de O''Higgins'
$ws.Range("S783").Value = 1000
$ws.Range("T783").Value = 15

